$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: sn1/sn2/sn3/sn4 -> fa1/fa2/fa3/fa4
$ws.Range("B1").Value = "fa1"
$ws.Range("C1").Value = "fa2"
$ws.Range("D1").Value = "fa3"
$ws.Range("E1").Value = "fa4"

# Highlight G4 with a yellow fill (same highlight used to flag the
# [M+Na]+ / [M+H]+ control bug fix mentioned in the commit message)
$ws.Range("G4").Interior.Color = 65535

# Move the active selection from G4 to E1
$ws.Range("E1").Select()
